$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2024-05-24 Friday", $false, $false, $false, $false, $false, $true, 1, $false, "2024-05-25 Saturday", 2) | Out-Null

# Update the answer cells in the table (row/col are 1-based)
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "23÷8=2, 7"  # was "52÷8=6, 4"
$t.Cell(1, 2).Range.Text = "28÷2=14, 0"  # was "43÷7=6, 1"
$t.Cell(1, 3).Range.Text = "20÷7=2, 6"  # was "60÷4=15, 0"
$t.Cell(1, 4).Range.Text = "57÷8=7, 1"  # was "63÷5=12, 3"
$t.Cell(1, 5).Range.Text = "97÷4=24, 1"  # was "76÷4=19, 0"

$t.Cell(5, 1).Range.Text = "17÷5=3, 2"  # was "38÷2=19, 0"
$t.Cell(5, 2).Range.Text = "29÷3=9, 2"  # was "40÷6=6, 4"
$t.Cell(5, 3).Range.Text = "53÷9=5, 8"  # was "22÷8=2, 6"
$t.Cell(5, 4).Range.Text = "76÷6=12, 4"  # was "93÷9=10, 3"
$t.Cell(5, 5).Range.Text = "76÷3=25, 1"  # was "84÷6=14, 0"

$t.Cell(9, 1).Range.Text = "30÷4=7, 2"  # was "76÷4=19, 0"
$t.Cell(9, 2).Range.Text = "73÷7=10, 3"  # was "95÷4=23, 3"
$t.Cell(9, 3).Range.Text = "15÷5=3, 0"  # was "33÷7=4, 5"
$t.Cell(9, 4).Range.Text = "27÷2=13, 1"  # was "40÷2=20, 0"
$t.Cell(9, 5).Range.Text = "47÷5=9, 2"  # was "56÷9=6, 2"

$t.Cell(13, 1).Range.Text = "53÷6=8, 5"  # was "32÷8=4, 0"
$t.Cell(13, 2).Range.Text = "43÷6=7, 1"  # was "72÷6=12, 0"
$t.Cell(13, 3).Range.Text = "47÷5=9, 2"  # was "17÷2=8, 1"
$t.Cell(13, 4).Range.Text = "12÷7=1, 5"  # was "58÷7=8, 2"
$t.Cell(13, 5).Range.Text = "59÷2=29, 1"  # was "81÷8=10, 1"

$t.Cell(17, 1).Range.Text = "27÷8=3, 3"  # was "87÷2=43, 1"
$t.Cell(17, 2).Range.Text = "53÷6=8, 5"  # was "61÷4=15, 1"
$t.Cell(17, 3).Range.Text = "73÷7=10, 3"  # was "42÷6=7, 0"
$t.Cell(17, 4).Range.Text = "57÷3=19, 0"  # was "98÷6=16, 2"
$t.Cell(17, 5).Range.Text = "33÷7=4, 5"  # was "15÷8=1, 7"
